# Update "想去人数" (want-to-go count) figures in column F for the
# 展览 (Exhibitions) and 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览": rows 3-7 map to F column updates ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 380
$ws1.Range("F4").Value = 4860
$ws1.Range("F5").Value = 20
$ws1.Range("F6").Value = 46
$ws1.Range("F7").Value = 485

# --- Sheet "全部类型": matching events, but shifted rows (6,8,9) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 380
$ws4.Range("F4").Value = 4860
$ws4.Range("F6").Value = 20
$ws4.Range("F8").Value = 46
$ws4.Range("F9").Value = 485
